# Summer2014 GE.xlsx - account for previously-missed classes (GEOG-450,
# GEOL-306, GEOL-311) and recompute OCNG-252 / WIEDERWOHL C percentages
# after fixing the course-delimiting bug described in the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Helper: write a plain text value into a cell while guaranteeing that
# Excel does NOT reinterpret look-alike numeric/percentage text as a
# number (which would change the cell type from shared-string to
# numeric) and without leaving a non-default cell style behind.
# ---------------------------------------------------------------------
function Set-TextValue {
    param($range, [string]$text)
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

# Everything from row 20 downward (the old GEOL-101 ... OCNG-252 block)
# is being rebuilt, so clear it first - this removes the stale rows and
# lets the sheet's used-range/dimension shrink-or-grow naturally.
$ws.Range("A20:H48").ClearContents()

# -----------------------------------------------------------------
# GEOG-450 (new course block)
# -----------------------------------------------------------------
Set-TextValue $ws.Range("A21") "GEOG-450"

Set-TextValue $ws.Range("B22") "PROUT E"
$ws.Range("C22").Value = 3.8
Set-TextValue $ws.Range("D22") "80.00%"
Set-TextValue $ws.Range("E22") "20.00%"
Set-TextValue $ws.Range("F22") "0.00%"
Set-TextValue $ws.Range("G22") "0.00%"
Set-TextValue $ws.Range("H22") "0.00%"

Set-TextValue $ws.Range("B23") "BRANNSTROM C"
$ws.Range("C23").Value = 3.857
Set-TextValue $ws.Range("D23") "85.71%"
Set-TextValue $ws.Range("E23") "14.29%"
Set-TextValue $ws.Range("F23") "0.00%"
Set-TextValue $ws.Range("G23") "0.00%"
Set-TextValue $ws.Range("H23") "0.00%"

# -----------------------------------------------------------------
# GEOL-101 (only MILLER B remains here now)
# -----------------------------------------------------------------
Set-TextValue $ws.Range("A25") "GEOL-101"

Set-TextValue $ws.Range("B26") "MILLER B"
$ws.Range("C26").Value = 3.1313
Set-TextValue $ws.Range("D26") "33.33%"
Set-TextValue $ws.Range("E26") "42.59%"
Set-TextValue $ws.Range("F26") "22.22%"
Set-TextValue $ws.Range("G26") "1.85%"
Set-TextValue $ws.Range("H26") "0.00%"

# -----------------------------------------------------------------
# GEOL-300
# -----------------------------------------------------------------
Set-TextValue $ws.Range("A28") "GEOL-300"

Set-TextValue $ws.Range("B29") "OLSZEWSKI T"
$ws.Range("C29").Value = 2.769
Set-TextValue $ws.Range("D29") "23.08%"
Set-TextValue $ws.Range("E29") "42.31%"
Set-TextValue $ws.Range("F29") "23.08%"
Set-TextValue $ws.Range("G29") "11.54%"
Set-TextValue $ws.Range("H29") "0.00%"

Set-TextValue $ws.Range("B30") "MILLER B"
$ws.Range("C30").Value = 2.769
Set-TextValue $ws.Range("D30") "11.54%"
Set-TextValue $ws.Range("E30") "65.38%"
Set-TextValue $ws.Range("F30") "11.54%"
Set-TextValue $ws.Range("G30") "11.54%"
Set-TextValue $ws.Range("H30") "0.00%"

Set-TextValue $ws.Range("B31") "HEANEY M"
$ws.Range("C31").Value = 3.28
Set-TextValue $ws.Range("D31") "28.00%"
Set-TextValue $ws.Range("E31") "72.00%"
Set-TextValue $ws.Range("F31") "0.00%"
Set-TextValue $ws.Range("G31") "0.00%"
Set-TextValue $ws.Range("H31") "0.00%"

# -----------------------------------------------------------------
# GEOL-306 (new course block)
# -----------------------------------------------------------------
Set-TextValue $ws.Range("A33") "GEOL-306"

Set-TextValue $ws.Range("B34") "TICE M"
$ws.Range("C34").Value = 4
Set-TextValue $ws.Range("D34") "100.00%"
Set-TextValue $ws.Range("E34") "0.00%"
Set-TextValue $ws.Range("F34") "0.00%"
Set-TextValue $ws.Range("G34") "0.00%"
Set-TextValue $ws.Range("H34") "0.00%"

# -----------------------------------------------------------------
# GEOL-311 (new course block)
# -----------------------------------------------------------------
Set-TextValue $ws.Range("A36") "GEOL-311"

Set-TextValue $ws.Range("B37") "MARCANTONIO F"
$ws.Range("C37").Value = 4
Set-TextValue $ws.Range("D37") "100.00%"
Set-TextValue $ws.Range("E37") "0.00%"
Set-TextValue $ws.Range("F37") "0.00%"
Set-TextValue $ws.Range("G37") "0.00%"
Set-TextValue $ws.Range("H37") "0.00%"

# -----------------------------------------------------------------
# GEOL-491 (TICE M belongs here, not under GEOL-491's old neighbour)
# -----------------------------------------------------------------
Set-TextValue $ws.Range("A39") "GEOL-491"

Set-TextValue $ws.Range("B40") "LAYA P"
$ws.Range("C40").Value = 4
Set-TextValue $ws.Range("D40") "100.00%"
Set-TextValue $ws.Range("E40") "0.00%"
Set-TextValue $ws.Range("F40") "0.00%"
Set-TextValue $ws.Range("G40") "0.00%"
Set-TextValue $ws.Range("H40") "0.00%"

Set-TextValue $ws.Range("B41") "GIARDINO J"
$ws.Range("C41").Value = 4
Set-TextValue $ws.Range("D41") "100.00%"
Set-TextValue $ws.Range("E41") "0.00%"
Set-TextValue $ws.Range("F41") "0.00%"
Set-TextValue $ws.Range("G41") "0.00%"
Set-TextValue $ws.Range("H41") "0.00%"

Set-TextValue $ws.Range("B42") "TICE M"
$ws.Range("C42").Value = 4
Set-TextValue $ws.Range("D42") "100.00%"
Set-TextValue $ws.Range("E42") "0.00%"
Set-TextValue $ws.Range("F42") "0.00%"
Set-TextValue $ws.Range("G42") "0.00%"
Set-TextValue $ws.Range("H42") "0.00%"

# -----------------------------------------------------------------
# GEOS-484
# -----------------------------------------------------------------
Set-TextValue $ws.Range("A44") "GEOS-484"

Set-TextValue $ws.Range("B45") "COLLINS D"
$ws.Range("C45").Value = 4
Set-TextValue $ws.Range("D45") "100.00%"
Set-TextValue $ws.Range("E45") "0.00%"
Set-TextValue $ws.Range("F45") "0.00%"
Set-TextValue $ws.Range("G45") "0.00%"
Set-TextValue $ws.Range("H45") "0.00%"

# -----------------------------------------------------------------
# OCNG-252 (GPA and percentage breakdown recomputed)
# -----------------------------------------------------------------
Set-TextValue $ws.Range("A47") "OCNG-252"

Set-TextValue $ws.Range("B48") "WIEDERWOHL C"
$ws.Range("C48").Value = 3.6712
Set-TextValue $ws.Range("D48") "70.42%"
Set-TextValue $ws.Range("E48") "25.35%"
Set-TextValue $ws.Range("F48") "4.23%"
Set-TextValue $ws.Range("G48") "0.00%"
Set-TextValue $ws.Range("H48") "0.00%"
